# Propuesta Compendio Estadístico Con Enfoque de Género 2022.docx
# "Cambios y/o correcciones en nombres de indicadores"

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $old"
    }
    return $ok
}

# 1) Shorten three indicator titles (drop the trailing ", según ..." clause)
Replace-Text "Mapa departamental por sexo, según departamento" "Mapa departamental por sexo"
Replace-Text "Mapa municipal por sexo, según municipio" "Mapa municipal por sexo"
Replace-Text "Esperanza de vida al nacer por sexo, según grupos de edad" "Esperanza de vida al nacer por sexo"

# 2) Drop stale <w:lastRenderedPageBreak/> markers on three indicator runs by
#    re-asserting the same text through Find/Replace, which rewrites the run
#    and clears the cached render marker.
Replace-Text "Tasa neta de escolaridad en el ciclo básico por sexo (serie histórica de 2018 a 2022)" "Tasa neta de escolaridad en el ciclo básico por sexo (serie histórica de 2018 a 2022)"
Replace-Text "Salario o ingresos promedio por sexo, según dominio de estudio y rama de actividad económica" "Salario o ingresos promedio por sexo, según dominio de estudio y rama de actividad económica"
Replace-Text "Mujeres magistradas en el Organismo Judicial " "Mujeres magistradas en el Organismo Judicial "

# 3) Fix "Tasa desempleo ..." -> "Tasa de desempleo ..."
Replace-Text "Tasa desempleo en la población de 15 años o más por sexo, según Pueblos " "Tasa de desempleo en la población de 15 años o más por sexo, según Pueblos "

# 4) Remove the whole "Tasa de matrimonios infantiles por sexo (serie histórica
#    de 2018 a 2022)" list item (duplicate/obsolete indicator).
$paras = $d.Paragraphs
$cnt = $paras.Count
for ($i = 1; $i -le $cnt; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "Tasa de matrimonios infantiles por sexo*") {
        $p.Range.Delete()
        break
    }
}
